# Update cryptocurrency price/volume data per the Jan 22 2023 symbol-list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "302.74"
Set-TextValue $ws.Range("E2") "-0.70%"

# Row 3
Set-TextValue $ws.Range("D3") "37.10"
Set-TextValue $ws.Range("E3") "6.21%"

# Row 4
Set-TextValue $ws.Range("D4") "5.021"
Set-TextValue $ws.Range("E4") "-3.13%"

# Row 5
Set-TextValue $ws.Range("D5") "0.07864"
Set-TextValue $ws.Range("E5") "0.47%"

# Row 6
Set-TextValue $ws.Range("D6") "2.212"
Set-TextValue $ws.Range("E6") "-4.35%"

# Row 7
Set-TextValue $ws.Range("E7") "-0.69%"

# Row 8
Set-TextValue $ws.Range("D8") "4.020"
Set-TextValue $ws.Range("E8") "0.75%"

# Row 9
Set-TextValue $ws.Range("D9") "0.9192"
Set-TextValue $ws.Range("E9") "-0.42%"

# Row 10
Set-TextValue $ws.Range("D10") "0.09557"
Set-TextValue $ws.Range("E10") "-5.39%"

# Row 11
Set-TextValue $ws.Range("D11") "0.1880"
Set-TextValue $ws.Range("E11") "2.24%"

# Row 12
Set-TextValue $ws.Range("D12") "0.08571"
Set-TextValue $ws.Range("E12") "0.04%"

# Row 13
Set-TextValue $ws.Range("D13") "0.03594"
Set-TextValue $ws.Range("E13") "6.54%"

# Row 14
Set-TextValue $ws.Range("D14") "0.09969"
Set-TextValue $ws.Range("E14") "0.64%"

# Row 15
Set-TextValue $ws.Range("D15") "0.001483"
Set-TextValue $ws.Range("E15") "0.21%"

# Row 16
Set-TextValue $ws.Range("D16") "0.005687"
Set-TextValue $ws.Range("E16") "0.89%"

# Row 17
Set-TextValue $ws.Range("E17") "-0.72%"

# Row 19
Set-TextValue $ws.Range("E19") "-0.66%"

# Row 20
Set-TextValue $ws.Range("E20") "-0.68%"

# Row 21
Set-TextValue $ws.Range("D21") "4.757"
Set-TextValue $ws.Range("E21") "4.27%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2200"
Set-TextValue $ws.Range("E22") "-7.86%"

# Row 23
Set-TextValue $ws.Range("D23") "0.04590"
Set-TextValue $ws.Range("E23") "-1.35%"

# Row 24
Set-TextValue $ws.Range("E24") "0.71%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004461"
Set-TextValue $ws.Range("E25") "0.20%"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001399"
Set-TextValue $ws.Range("E26") "7.71%"

# Row 27
Set-TextValue $ws.Range("E27") "39.83%"

# Row 39
Set-TextValue $ws.Range("D39") "0.01804"
Set-TextValue $ws.Range("E39") "3.67%"

# Row 40
Set-TextValue $ws.Range("D40") "0.04717"
Set-TextValue $ws.Range("E40") "-0.80%"

# Row 41
Set-TextValue $ws.Range("D41") "0.008126"
Set-TextValue $ws.Range("E41") "5.65%"

# Row 42
Set-TextValue $ws.Range("E42") "-1.19%"

# Row 43
Set-TextValue $ws.Range("D43") "0.007551"
Set-TextValue $ws.Range("E43") "6.81%"

# Row 44
Set-TextValue $ws.Range("D44") "0.002229"
Set-TextValue $ws.Range("E44") "0.96%"

# Row 45
Set-TextValue $ws.Range("D45") "0.01041"
Set-TextValue $ws.Range("E45") "1.88%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006160"
Set-TextValue $ws.Range("E46") "2.77%"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.02%"

# Row 48
Set-TextValue $ws.Range("D48") "0.0005801"
Set-TextValue $ws.Range("E48") "0.01%"

# Row 49
Set-TextValue $ws.Range("D49") "7.087"
Set-TextValue $ws.Range("E49") "22.26%"

# Row 50
Set-TextValue $ws.Range("E50") "0.09%"

# Row 51
Set-TextValue $ws.Range("D51") "0.00002099"
Set-TextValue $ws.Range("E51") "0.02%"
